$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row for the new "High Voltage (1200V)" columns ---
$ws.Range("J1").Value = "High Voltage (1200V)"
$ws.Range("J2").Value = "1120V"
$ws.Range("I1").Value = "Expected(1200V)"

# Match the header styling (font/alignment) used by the rest of row 1,
# copied from existing cells so no new style entries are created.
$ws.Range("B1:C1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Column I: Expected (1200V) = E * 120 ---
$ws.Range("I2").Formula = "=E2*120"
$ws.Range("I3:I10").Formula = "=E3*120"

# --- Column J: measured high-voltage (1200V) values (entered out of row
# order, matching the original author's data-entry sequence so the shared
# string table comes out in the same order) ---
$ws.Range("J6").Value = "581.5 (changing)"
$ws.Range("J3").Value = "836 (changing)"
$ws.Range("J4").Value = "740 (changing)"
$ws.Range("J5").Value = "640 (changing)"
$ws.Range("J7").Value = "508.6 (changing)"
$ws.Range("J8").Value = "443.6 (changing)"
$ws.Range("J9").Value = "381 (changing)"
$ws.Range("J10").Value = "319 (changing)"

# J8:J10 pick up the right-aligned numeric style used by column G, copied
# from the existing G column cells (reuses the existing style, no new
# style entries added).
$ws.Range("G8:G10").Copy()
$ws.Range("J8:J10").PasteSpecial(-4122)

# New column width for column I
$ws.Columns.Item(9).ColumnWidth = 29.15

# Selection moves to J10 after the edit
$ws.Range("J10").Select()
